# Add a "Theory" hyperlink after the Normalization bullet on the
# "Data modelling" slide (slide 8), per commit:
#   "Added link to database normalization markdown from PPT"
#
# Before: "Normalization – Minimize redundancy and dependency (1NF, 2NF, 3NF)"
# After : "Normalization – Minimize redundancy and dependency (1NF, 2NF, 3NF) - Theory"
#          (" - Theory" appended; " " and "Theory" are coloured blue and
#           "Theory" carries the hyperlink.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shp = $s.Shapes.Item("Content Placeholder 2")

$tr = $shp.TextFrame.TextRange

# Find the paragraph that starts the Normalization bullet.
$normPara = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text.StartsWith("Normalization")) {
        $normPara = $para
    }
}

# TextRange.Text includes the trailing paragraph-mark character; strip it
# off before appending so the new text extends this paragraph instead of
# spilling into a new one.
$paraLen = $normPara.Length
$coreText = $normPara.Characters(1, $paraLen - 1).Text

# Append " - Theory" -> the " -" part keeps the original (unlinked) run
# formatting; the " " and "Theory" that follow are recoloured/linked below.
$suffix = " - Theory"
$normPara.Text = $coreText + $suffix
$newLen = $coreText.Length + $suffix.Length

# Blue (RGB 0,112,192 = hex 0070C0), matching the link colour used for
# hyperlinks in this deck.
$linkColor = 12611584

# The single space before "Theory".
$spaceRun = $normPara.Characters($newLen - 6, 1)
$spaceRun.Font.Color.RGB = $linkColor

# "Theory" itself: colour + hyperlink to the normalization write-up.
$theoryRun = $normPara.Characters($newLen - 5, 6)
$theoryRun.Font.Color.RGB = $linkColor
$theoryRun.ActionSettings(1).Hyperlink.Address = "Normalization.md"
